$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the "TextBox 83" shape that contains "Possible values for each feature"
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 83") {
        $shp.Delete()
    }
}
